$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 39 and 40: Monero and ImmutableX swap places (with updated price/volume data)
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.58"
$ws.Range("E39").Value = "  +7.22%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'170.13"
$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("D2").Value = "63.384.85"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "3.479.01"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'580.77"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'147.75"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("D7").Value = "3.480.79"
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'7.72"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").Value = "'0.406"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").Value = "4.068.68"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "'29.70"
$ws.Range("E14").Value = "  +6.11%  "
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").Value = "3.470.51"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "'0.0000172"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "63.303.88"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "'6.35"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("D20").Value = "'14.45"
$ws.Range("E20").Value = "  +4.90%  "
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").Value = "'390.42"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").Value = "'75.00"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "3.607.32"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").Value = "'7.61"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'8.23"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "'7.15"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "'5.31"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").Value = "'32.21"
$ws.Range("E38").Value = "  +14.00%  "
$ws.Range("D41").Value = "3.509.72"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "'0.0765"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "'0.800"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  +4.77%  "
$ws.Range("D46").Value = "'1.73"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").Value = "'4.42"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "2.629.89"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").Value = "'2.32"
$ws.Range("E49").Value = "  +11.36%  "
$ws.Range("D50").Value = "'23.19"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "'6.79"
$ws.Range("E51").Value = "  +2.27%  "
